$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.824475333333333
$ws.Range("N2").Value = 5.473426
$ws.Range("O2").Value = 0.1906606574278047
$ws.Range("P2").Value = 0.2015451970524477
$ws.Range("Q2").Value = 0.7086286521004445
$ws.Range("R2").Value = 6.377657868904
$ws.Range("S2").Value = 0.1906606574278047
$ws.Range("T2").Value = 0.2015451970524477

$ws.Range("O3").Value = 0.6423822165107047
$ws.Range("P3").Value = 0.6790548829333741
$ws.Range("S3").Value = 0.6423822165107047
$ws.Range("T3").Value = 0.6790548829333741

$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009389666666666666
$ws.Range("N4").Value = 0.028169
$ws.Range("O4").Value = 0.0009812355294625031
$ws.Range("P4").Value = 0.001037252838673693
$ws.Range("Q4").Value = 0.003646959052888889
$ws.Range("R4").Value = 0.032822631476
$ws.Range("S4").Value = 0.0009812355294625031
$ws.Range("T4").Value = 0.001037252838673693

$ws.Range("M5").Value = 1.5503715
$ws.Range("N5").Value = 3.100743
$ws.Range("O5").Value = 0.1620163583726162
$ws.Range("P5").Value = 0.1141770910840848
$ws.Range("Q5").Value = 0.6021663577620001
$ws.Range("R5").Value = 3.612998146572
$ws.Range("S5").Value = 0.1620163583726162
$ws.Range("T5").Value = 0.1141770910840848

$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03788966666666667
$ws.Range("N6").Value = 0.113669
$ws.Range("O6").Value = 0.003959532159411881
$ws.Range("P6").Value = 0.004185576091419648
$ws.Range("Q6").Value = 0.01471639705288889
$ws.Range("R6").Value = 0.132447573476
$ws.Range("S6").Value = 0.003959532159411881
$ws.Range("T6").Value = 0.004185576091419648
